# Middleware test results - insert original method results
# (see commit message: "updated excel file with the insert original method results")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H ("insert" / Original method, "In Seconds") data updates ---
# Existing rows 35-41 get new timing values (the old values effectively
# shift down as new measurements were inserted ahead of them), and rows
# 42-44 (previously blank) now get values too.
$ws.Range("H35").Value = 127.914
$ws.Range("H36").Value = 128.842
$ws.Range("H37").Value = 120.159
$ws.Range("H38").Value = 129.10399999999998
$ws.Range("H39").Value = 128.958
$ws.Range("H40").Value = 129.535
$ws.Range("H41").Value = 120.803
$ws.Range("H42").Value = 122.167
$ws.Range("H43").Value = 123.516
$ws.Range("H44").Value = 124.818

# --- Remove the yellow highlight fill from the "insert" header (H4) and
#     from the column-H average cell (H45); both previously used a
#     yellow solid fill and now use no fill ("No Fill" / xlPatternNone). ---
$ws.Range("H4").Interior.Pattern = -4142
$ws.Range("H45").Interior.Pattern = -4142

# --- Selection / view state: the active selection moved to the newly
#     populated "insert" results block (L5:R46) with the view scrolled
#     so row 5 is at the top. ---
$ws.Range("L5:R46").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
